# Scheduled market-price refresh for the Unicorn_Profits leve-profit tables.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) on
# several per-job sheets to reflect newly pulled marketboard data.
$wb = $excel.ActiveWorkbook

# --- ALC: "Asking for a Friend" (row 100) price refresh ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5199.4517
$ws.Range("I100").Value = 1386.875
$ws.Range("J100").Value = 6525.5654
$ws.Range("K100").Value = 1386.875
$ws.Range("L100").Value = 6525.5654
$ws.Range("M100").Value = -845.875
$ws.Range("N100").Value = -7607.5654

# --- ARM: rows 121-141 no longer have marketboard data; clear H:N ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:N127").ClearContents()
$ws.Range("H128:N128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

# --- BSM: "The Gold Experience" (row 107) price refresh ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2056
$ws.Range("I107").Value = 1405.5
$ws.Range("J107").Value = 2706.5
$ws.Range("K107").Value = 1405.5
$ws.Range("L107").Value = 2706.5
$ws.Range("M107").Value = 514.5
$ws.Range("N107").Value = -6546.5

# --- CRP: rows 129-141 now have newly-pulled marketboard data (H:N) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 38500
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 38500
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 38500
$ws.Range("N129").Value = -48500
$ws.Range("H130").Value = 34653.332
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 34653.332
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 34653.332
$ws.Range("N130").Value = -44693.332
$ws.Range("H131").Value = 21000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 21000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 21000
$ws.Range("N131").Value = -31080
$ws.Range("H132").Value = 1784.9131
$ws.Range("I132").Value = 925.8788
$ws.Range("J132").Value = 3965.5386
$ws.Range("K132").Value = 2777.6364
$ws.Range("L132").Value = 11896.6158
$ws.Range("M132").Value = -247.6363999999999
$ws.Range("N132").Value = -16956.6158
$ws.Range("H133").Value = 40975.2
$ws.Range("I133").Value = 21296
$ws.Range("J133").Value = 45895
$ws.Range("K133").Value = 21296
$ws.Range("L133").Value = 45895
$ws.Range("M133").Value = -18766
$ws.Range("N133").Value = -50955
$ws.Range("H134").Value = 1532.625
$ws.Range("I134").Value = 1023.32
$ws.Range("J134").Value = 2381.4666
$ws.Range("K134").Value = 3069.96
$ws.Range("L134").Value = 7144.399800000001
$ws.Range("M134").Value = -534.96
$ws.Range("N134").Value = -12214.3998
$ws.Range("H135").Value = 69487.5
$ws.Range("I135").Value = 69487.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 69487.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -64417.5
$ws.Range("H137").Value = 41320
$ws.Range("I137").Value = 16500
$ws.Range("J137").Value = 49593.332
$ws.Range("K137").Value = 16500
$ws.Range("L137").Value = 49593.332
$ws.Range("M137").Value = -11400
$ws.Range("N137").Value = -59793.332
$ws.Range("H138").Value = 33520
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 33520
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 33520
$ws.Range("N138").Value = -43800
$ws.Range("H139").Value = 39378
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 39378
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 39378
$ws.Range("N139").Value = -49658
$ws.Range("H140").Value = 32945
$ws.Range("I140").Value = 20000
$ws.Range("J140").Value = 37260
$ws.Range("K140").Value = 20000
$ws.Range("L140").Value = 37260
$ws.Range("M140").Value = -14820
$ws.Range("N140").Value = -47620
$ws.Range("H141").Value = 31572.25
$ws.Range("I141").Value = 25300
$ws.Range("J141").Value = 33663
$ws.Range("K141").Value = 25300
$ws.Range("L141").Value = 33663
$ws.Range("M141").Value = -20120
$ws.Range("N141").Value = -44023

# --- CUL: row 122 price refresh ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 748.8182
$ws.Range("I122").Value = 322.375
$ws.Range("J122").Value = 1150.1765
$ws.Range("K122").Value = 2901.375
$ws.Range("L122").Value = 10351.5885
$ws.Range("M122").Value = -451.375
$ws.Range("N122").Value = -15251.5885

# --- GSM: rows 88, 91, 113, 126 price refresh ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 39000
$ws.Range("J88").Value = 39000
$ws.Range("L88").Value = 39000
$ws.Range("N88").Value = -39902
$ws.Range("H91").Value = 39000
$ws.Range("J91").Value = 39000
$ws.Range("L91").Value = 39000
$ws.Range("N91").Value = -42120
$ws.Range("H113").Value = 1505.2858
$ws.Range("I113").Value = 791.8333
$ws.Range("J113").Value = 2040.375
$ws.Range("K113").Value = 791.8333
$ws.Range("L113").Value = 2040.375
$ws.Range("M113").Value = 1378.1667
$ws.Range("N113").Value = -6380.375
$ws.Range("H126").Value = 4710.5625
$ws.Range("I126").Value = 4387.5
$ws.Range("J126").Value = 5679.75
$ws.Range("K126").Value = 13162.5
$ws.Range("L126").Value = 17039.25
$ws.Range("M126").Value = -10692.5
$ws.Range("N126").Value = -21979.25
